$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells: "_old" suffix -> "_FV2210", "_new" suffix -> "_FV2304" ---
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value2
    if ($v -ne $null) {
        if ($v.EndsWith("_old")) {
            $cell.Value2 = $v.Substring(0, $v.Length - 4) + "_FV2210"
        } elseif ($v.EndsWith("_new")) {
            $cell.Value2 = $v.Substring(0, $v.Length - 4) + "_FV2304"
        }
    }
}

# --- 2. Convert the data range A1:U55 into an Excel Table ("Table1") with an AutoFilter ---
# Stash the header row's existing formatting (bold/fill/border/wrap) on a scratch
# cell outside the used range, clear the header's direct formatting so
# ListObjects.Add doesn't need to capture a header-style override (dxf), then
# restore the original formatting in a single paste-formats operation so no
# extra/orphaned style records are left behind in styles.xml.
$hdr = $ws.Range("A1:U1")
$scratch = $ws.Range("Z1")
$ws.Range("A1").Copy() | Out-Null
$scratch.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$hdr.ClearFormats()

$range = $ws.Range("A1:U55")
$tbl = $ws.ListObjects.Add(1, $range, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

$scratch.Copy() | Out-Null
$hdr.PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$excel.CutCopyMode = $false
$scratch.Clear()

# --- 3. Freeze the header row (split after row 1) ---
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
